$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Repair the pre-existing ISO "t=d" date cells (A5:A36) so they hold a
#    proper date value instead of an error. Re-assigning a date-like string
#    to a cell that already carries the date style (s="1") keeps that same
#    style (numFmtId 14) and simply stores the correct serial number.
# ---------------------------------------------------------------------------
$dates = @{
    5  = "10/04/2022"
    6  = "10/07/2022"
    7  = "10/07/2022"
    8  = "10/10/2022"
    9  = "10/10/2022"
    10 = "10/11/2022"
    11 = "10/14/2022"
    12 = "10/14/2022"
    13 = "10/20/2022"
    14 = "10/20/2022"
    15 = "10/25/2022"
    16 = "10/25/2022"
    18 = "10/28/2022"
    19 = "10/28/2022"
    20 = "10/30/2022"
    21 = "10/31/2022"
    22 = "10/31/2022"
    23 = "11/04/2022"
    24 = "11/05/2022"
    25 = "11/07/2022"
    26 = "11/08/2022"
    27 = "12/01/2022"
    28 = "12/02/2022"
    29 = "12/02/2022"
    30 = "12/09/2022"
    31 = "12/09/2022"
    32 = "12/12/2022"
    33 = "12/13/2022"
    34 = "12/18/2022"
    35 = "12/30/2022"
    36 = "01/04/2023"
}
foreach ($row in $dates.Keys) {
    $ws.Cells.Item($row, 1).Value = $dates[$row]
}

# ---------------------------------------------------------------------------
# 2) Row 17: "Created Bmp File" -> "Created Bmp converter", and merge C17:E17
#    like the other task cells (was a lone cell before).
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = "Created Bmp converter"
$ws.Range("C17:E17").Merge()

# ---------------------------------------------------------------------------
# 3) Rows 34-36 get their C:E task cell merged too (text stays the same).
# ---------------------------------------------------------------------------
$ws.Range("C34:E34").Merge()
$ws.Range("C35:E35").Merge()
$ws.Range("C36:E36").Merge()

# ---------------------------------------------------------------------------
# 4) Append new row 37: D. Hoyer, 2023-01-04, "Score funktioniert und jump
#    angepasst", 2.5h.
# ---------------------------------------------------------------------------
$ws.Range("A36").Copy()
$ws.Range("A37").PasteSpecial(-4122)   # xlPasteFormats - reuse the date style
$ws.Range("A37").Value = "01/04/2023"

$ws.Range("B37").Value = "D. Hoyer"

$ws.Range("C37").Value = "Score funktioniert und jump angepasst"
$ws.Range("C37:E37").Merge()

$ws.Range("F37").Value = 2.5
$ws.Range("G37").Value = "h"

# ---------------------------------------------------------------------------
# 5) Selection moves to F38 (first empty row under the new data) - this also
#    clears the stale topLeftCell="A21" scroll position.
# ---------------------------------------------------------------------------
$ws.Range("F38").Select()

# ---------------------------------------------------------------------------
# 6) Recalculate so the SUMIF/SUM totals in row 4 pick up the new row.
# ---------------------------------------------------------------------------
$excel.CalculateFullRebuild()

Write-Host "edit complete"
